$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 243, shifting existing rows 243:267 down to 244:268.
$ws.Rows.Item(243).Insert()

# Copy static / descriptive columns from the row below (now row 244) which still
# holds the original "row 243" values for these shared columns.
$ws.Cells.Item(243, 1).Value2 = $ws.Cells.Item(244, 1).Value2   # A Mercado ID
$ws.Cells.Item(243, 2).Value2 = $ws.Cells.Item(244, 2).Value2   # B Mercado
$ws.Cells.Item(243, 3).Value2 = $ws.Cells.Item(244, 3).Value2   # C Region
$ws.Cells.Item(243, 5).Value2 = $ws.Cells.Item(244, 5).Value2   # E Codreg
$ws.Cells.Item(243, 6).Value2 = $ws.Cells.Item(244, 6).Value2   # F Tipo
$ws.Cells.Item(243, 7).Value2 = $ws.Cells.Item(244, 7).Value2   # G Producto ID
$ws.Cells.Item(243, 8).Value2 = $ws.Cells.Item(244, 8).Value2   # H Producto
$ws.Cells.Item(243, 9).Value2 = $ws.Cells.Item(244, 9).Value2   # I Categoria ID
$ws.Cells.Item(243, 10).Value2 = $ws.Cells.Item(244, 10).Value2 # J Categoria

# New record values for the inserted row.
$ws.Cells.Item(243, 4).Value2 = 44776                 # D Fecha
$ws.Cells.Item(243, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(243, 11).Value2 = "Sin especificar"     # K Variedad
$ws.Cells.Item(243, 12).Value2 = "2a amarillo"         # L Calidad
$ws.Cells.Item(243, 13).Value2 = 250                   # M Volumen
$ws.Cells.Item(243, 14).Value2 = 9000                  # N Precio minimo
$ws.Cells.Item(243, 15).Value2 = 10000                 # O Precio maximo
$ws.Cells.Item(243, 16).Value2 = 9500                  # P Precio promedio ponderado
$ws.Cells.Item(243, 17).Value2 = "`$/caja 20 kilos"     # Q Unidad de comercializacion
$ws.Cells.Item(243, 18).Value2 = "Región de Coquimbo"   # R Origen
$ws.Cells.Item(243, 19).Value2 = 475                   # S Precio $/Kg
$ws.Cells.Item(243, 20).Value2 = 20                    # T Kg / unidad

$wb.Save()
